$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos (values refresh, label unchanged) ---
$ws.Range("B4").Value = 1364447
$ws.Range("C4").Value = 17138
$ws.Range("E4").Value = 1043141
$ws.Range("F4").Value = 16490
$ws.Range("G4").Value = 653
$ws.Range("H4").Value = 80690

# --- Row 11: Brasil (values refresh, label unchanged) ---
$ws.Range("B11").Value = 157695
$ws.Range("C11").Value = 1634
$ws.Range("E11").Value = 85253
$ws.Range("G11").Value = 101
$ws.Range("H11").Value = 10757

# --- Rows 16/17: India & Peru swap order + data refresh ---
# Row 16 becomes Peru (fresh data), Row 17 becomes India (old row16 data)
$ws.Range("A16").Value = "Peru"
$ws.Range("B16").Value = 67307
$ws.Range("C16").Value = 2292
$ws.Range("D16").Value = 20246
$ws.Range("E16").Value = 45172
$ws.Range("F16").Value = 774
$ws.Range("G16").Value = 75
$ws.Range("H16").Value = 1889

$ws.Range("A17").Value = "India"
$ws.Range("B17").Value = 67161
$ws.Range("C17").Value = 4353
$ws.Range("D17").Value = 20969
$ws.Range("E17").Value = 43980
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 111
$ws.Range("H17").Value = 2212

# --- Rows 77/78: Bosnia y Herzegovina & Guinea swap order + data refresh ---
# Row 77 becomes Guinea (fresh data), Row 78 becomes Bosnia y Herzegovina (old row77 data)
$ws.Range("A77").Value = "Guinea"
$ws.Range("B77").Value = 2146
$ws.Range("C77").Value = 104
$ws.Range("D77").Value = 714
$ws.Range("E77").Value = 1421
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 11

$ws.Range("A78").Value = "Bosnia y Herzegovina"
$ws.Range("B78").Value = 2117
$ws.Range("C78").Value = 27
$ws.Range("D78").Value = 1106
$ws.Range("E78").Value = 904
$ws.Range("F78").Value = 4
$ws.Range("G78").Value = 5
$ws.Range("H78").Value = 107

# --- Row 85: Costa de Marfil (values refresh, label unchanged) ---
$ws.Range("B85").Value = 1700
$ws.Range("C85").Value = 33
$ws.Range("D85").Value = 794
$ws.Range("E85").Value = 885

# --- Row 107: Niger (values refresh, label unchanged) ---
$ws.Range("B107").Value = 821
$ws.Range("C107").Value = 6
$ws.Range("D107").Value = 624
$ws.Range("E107").Value = 151
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 46

# --- Rows 112-115: Paraguay moves from after Uruguay to before Crucero ---
# Crucero, Mali, Uruguay each shift down one row; Paraguay gets fresh data at row 112
$ws.Range("A112").Value = "Paraguay"
$ws.Range("B112").Value = 713
$ws.Range("C112").Value = 24
$ws.Range("D112").Value = 165
$ws.Range("E112").Value = 538
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 10

$ws.Range("A113").Value = "Crucero"
$ws.Range("B113").Value = 712
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 645
$ws.Range("E113").Value = 54
$ws.Range("F113").Value = 4
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 13

$ws.Range("A114").Value = "Mali"
$ws.Range("B114").Value = 704
$ws.Range("C114").Value = 12
$ws.Range("D114").Value = 351
$ws.Range("E114").Value = 315
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 38

$ws.Range("A115").Value = "Uruguay"
$ws.Range("B115").Value = 702
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 513
$ws.Range("E115").Value = 171
$ws.Range("F115").Value = 8
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 18

Write-Host "Applied all updates"
